# Split the "Consumo" paragraph's single run into two runs:
#   1) "El consumo no ha de superar 1 A de corriente"
#   2) ", para evitar un sobrecalentamiento del dispositivo que pudiese
#      tener fatales consecuencias."
# Both runs keep the original formatting (sz=24 / szCs=24).

$d = $word.ActiveDocument

# Locate the full original sentence (including the trailing ". ").
$target = $d.Content
$found = $target.Find.Execute(
    "El consumo no ha de superar 1 A de corriente. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Consumo' sentence to edit."
}

$sentenceStart = $target.Start
$sentenceEnd   = $target.End

# Split point: right after "...de corriente" (before the ". " tail).
$splitPoint = $sentenceStart + 44

# First half keeps its original wording untouched.
$firstRange = $d.Range($sentenceStart, $splitPoint)

# Second half (currently ". ") becomes the new clause, replacing the
# old trailing period + space.
$secondRange = $d.Range($splitPoint, $sentenceEnd)
$newClause = ", para evitar un sobrecalentamiento del dispositivo que pudiese tener fatales consecuencias."
$secondRange.Text = $newClause

# Re-grab the range that now holds the inserted clause so we can nudge
# its formatting; toggling a character property and restoring it keeps
# this text in its own run instead of being silently re-merged with the
# preceding, identically-formatted run when the document is saved.
$secondRangeNow = $d.Range($splitPoint, $splitPoint + $newClause.Length)
$secondRangeNow.Font.Bold = $true
$secondRangeNow.Font.Bold = $false

Write-Output "Split 'Consumo' sentence into two runs."
